$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text cell (matches the source data, which
    # stores every Price/Volume figure as a string, incl. plain-numeric-looking ones).
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row identity swaps: rows 44/45 and 47/48 exchanged their coin data ---
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D44") "0.785"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D45") "1.00"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D47") "1.96"
$ws.Range("E47").Value = "  +2.56%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D48") "159.04"
$ws.Range("E48").Value = "  -1.10%  "

# --- Price / Volume(1h) refresh for all other rows ---
$ws.Range("D2").Value = "98.975.63"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.297.78"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "255.18"
$ws.Range("E5").Value = "  +0.54%  "
Set-TextValue $ws.Range("D6") "627.88"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +23.20%  "
Set-TextValue $ws.Range("D8") "0.404"
$ws.Range("E8").Value = "  +5.34%  "
Set-TextValue $ws.Range("D10") "0.982"
$ws.Range("E10").Value = "  +23.26%  "
$ws.Range("D11").Value = "3.294.24"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("E12").Value = "  +2.83%  "
Set-TextValue $ws.Range("D13") "40.87"
$ws.Range("E13").Value = "  +14.56%  "
$ws.Range("D14").Value = "98.686.78"
$ws.Range("E14").Value = "  +1.19%  "
Set-TextValue $ws.Range("D15") "0.0000251"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("D16").Value = "3.913.41"
$ws.Range("E16").Value = "  -0.84%  "
Set-TextValue $ws.Range("D17") "5.49"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "3.304.16"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  -3.47%  "
Set-TextValue $ws.Range("D20") "15.67"
$ws.Range("E20").Value = "  +6.54%  "
Set-TextValue $ws.Range("D21") "6.42"
$ws.Range("E21").Value = "  +9.59%  "
Set-TextValue $ws.Range("D22") "488.64"
$ws.Range("E22").Value = "  +2.02%  "
Set-TextValue $ws.Range("D23") "9.44"
$ws.Range("E23").Value = "  +4.06%  "
Set-TextValue $ws.Range("D24") "0.0000204"
$ws.Range("E24").Value = "  -1.00%  "
Set-TextValue $ws.Range("D25") "5.73"
$ws.Range("E25").Value = "  +1.20%  "
Set-TextValue $ws.Range("D26") "0.341"
$ws.Range("E26").Value = "  +38.09%  "
Set-TextValue $ws.Range("D27") "90.47"
$ws.Range("E27").Value = "  +3.42%  "
Set-TextValue $ws.Range("D28") "12.19"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("D29").Value = "3.463.97"
$ws.Range("E29").Value = "  -3.00%  "
Set-TextValue $ws.Range("D30") "0.148"
$ws.Range("E30").Value = "  +20.58%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +1.74%  "
$ws.Range("E33").Value = "  +17.53%  "
Set-TextValue $ws.Range("D34") "0.998"
$ws.Range("E34").Value = "  -0.21%  "
Set-TextValue $ws.Range("D35") "28.04"
$ws.Range("E35").Value = "  +3.55%  "
Set-TextValue $ws.Range("D36") "0.484"
$ws.Range("E36").Value = "  +8.52%  "
$ws.Range("E37").Value = "  +0.45%  "
Set-TextValue $ws.Range("D38") "7.34"
$ws.Range("E38").Value = "  +1.58%  "
$ws.Range("E39").Value = "  +1.46%  "
Set-TextValue $ws.Range("D40") "497.03"
$ws.Range("E40").Value = "  -4.38%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("E43").Value = "  -0.72%  "
Set-TextValue $ws.Range("D46") "3.17"
$ws.Range("E46").Value = "  -1.32%  "
Set-TextValue $ws.Range("D49") "4.85"
$ws.Range("E49").Value = "  +8.21%  "
$ws.Range("E50").Value = "  +16.34%  "
Set-TextValue $ws.Range("D51") "0.850"
$ws.Range("E51").Value = "  +7.75%  "
